# issue #5: stock data from json to db
# Add "category", "source_file", "index" columns to the 股票 (stock) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 20

# Insert a new column I ("category") between the existing "property_category"
# (H) and "date" (old I, now shifted to J) columns. Insert() shifts everything
# from I onward one column to the right (I->J, J->K, K->L).
$ws.Columns.Item(9).Insert()

# Header + values for the newly inserted "category" column.
$ws.Range("I1").Value = "category"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
}

# Append two new trailing columns: "source_file" and "index".
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "tmp84bd1"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}
